$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 92, column A (date/time) value
$ws.Range("A92").Value = 45455.2916666667

# Add new row 93 with data (results from R script)
$ws.Range("A93").Value = 45456.2916666667
$ws.Range("B93").Value = 11300
$ws.Range("C93").Value = 0.735000014305115
$ws.Range("D93").Value = 0.714999973773956
$ws.Range("E93").Value = 0.735000014305115
$ws.Range("F93").Value = 0.714999973773956

# adj_close (G) is stored as text in this sheet, same as the rest of the column
$ws.Range("G93").Value = "'0.714999973773956"
$ws.Range("G93").Style = "Normal"

$ws.Range("H93").Value = "BWZ.MI"

# Copy the date formatting from A92 onto the new A93 cell
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$excel.CutCopyMode = $false
